$wb = $excel.ActiveWorkbook
$hh = $wb.Worksheets.Item("HH")
$presupuesto = $wb.Worksheets.Item("Presupuesto")

# --- HH sheet ---

# S3: bump the hourly rate used in the formula (284.03 -> 314.25)
$hh.Range("S3").Formula = "=(314.25)*S2"

# New R12 entry (additional hours logged), matching the time-duration
# formatting already used by its neighbours R8:R11
$hh.Range("R12").Value = 1.2652777777777777
$hh.Range("R12").NumberFormat = "[h]:mm:ss"

# Q6: extend the sum to include the newly tracked R12 hours (must be set
# after R12 has a value so the cached result recalculates correctly)
$hh.Range("Q6").Formula = "= Q8+R8+R9+R10+R11+R12"

# Q8 switches to a date/time number format (adds a new cell style)
$hh.Range("Q8").NumberFormat = "m/d/yy h:mm"

# Column Q widens to fit the new date/time content
$hh.Columns.Item(17).ColumnWidth = 13.5

# Move the active selection to S4
$hh.Range("S4").Select()

# --- Presupuesto sheet ---
$presupuesto.Activate()

# Move the active selection to F10
$presupuesto.Range("F10").Select()
